# Auto-generated script applying odds updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("G5").Value = 2.35
$ws.Range("I5").Value = 2.75
$ws.Range("U5").Value = 13
$ws.Range("X5").Value = 19
$ws.Range("AE5").Value = 10
# Row 8
$ws.Range("G8").Value = 1.42
$ws.Range("I8").Value = 7.1
$ws.Range("M8").Value = 3.25
$ws.Range("O8").Value = 1.88
$ws.Range("R8").Value = 1.93
$ws.Range("U8").Value = 6.4
$ws.Range("W8").Value = 9.25
$ws.Range("Z8").Value = 10.75
$ws.Range("AB8").Value = 19.5
$ws.Range("AC8").Value = 100
$ws.Range("AD8").Value = 800
$ws.Range("AE8").Value = 18
$ws.Range("AG8").Value = 22
$ws.Range("AH8").Value = 175
$ws.Range("AI8").Value = 80
# Row 11
$ws.Range("G11").Value = 2.47
$ws.Range("H11").Value = 2.32
$ws.Range("J11").Value = 1.23
$ws.Range("K11").Value = 3.7
$ws.Range("L11").Value = 1.87
$ws.Range("M11").Value = 1.83
$ws.Range("N11").Value = 3.55
$ws.Range("O11").Value = 1.26
$ws.Range("P11").Value = 1.82
$ws.Range("Q11").Value = 1.88
$ws.Range("R11").Value = 2.52
$ws.Range("S11").Value = 1.47
$ws.Range("T11").Value = 4.75
$ws.Range("V11").Value = 11.25
$ws.Range("Y11").Value = 70
$ws.Range("Z11").Value = 3.7
$ws.Range("AA11").Value = 5.2
$ws.Range("AB11").Value = 23
$ws.Range("AC11").Value = 200
$ws.Range("AE11").Value = 6.8
$ws.Range("AG11").Value = 16
$ws.Range("AH11").Value = 80
$ws.Range("AI11").Value = 65
$ws.Range("AJ11").Value = 100
# Row 12
$ws.Range("G12").Value = 3.6
$ws.Range("H12").Value = 2.45
$ws.Range("I12").Value = 2.55
$ws.Range("J12").Value = 1.17
$ws.Range("K12").Value = 4.45
$ws.Range("L12").Value = 1.62
$ws.Range("M12").Value = 2.15
$ws.Range("N12").Value = 2.85
$ws.Range("O12").Value = 1.37
$ws.Range("P12").Value = 1.65
$ws.Range("Q12").Value = 2.12
$ws.Range("R12").Value = 2.12
$ws.Range("S12").Value = 1.65
$ws.Range("T12").Value = 7.2
$ws.Range("V12").Value = 13
$ws.Range("W12").Value = 60
$ws.Range("X12").Value = 45
$ws.Range("Y12").Value = 60
$ws.Range("Z12").Value = 4.45
$ws.Range("AA12").Value = 5
$ws.Range("AB12").Value = 17
$ws.Range("AC12").Value = 120
$ws.Range("AE12").Value = 5.7
$ws.Range("AF12").Value = 11.25
$ws.Range("AG12").Value = 10
$ws.Range("AH12").Value = 30
$ws.Range("AI12").Value = 28
$ws.Range("AJ12").Value = 50
# Row 14
$ws.Range("G14").Value = 2.3
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 3.2
$ws.Range("K14").Value = 5.6
$ws.Range("L14").Value = 1.5
$ws.Range("M14").Value = 2.42
$ws.Range("N14").Value = 2.42
$ws.Range("P14").Value = 1.52
$ws.Range("R14").Value = 2.05
$ws.Range("S14").Value = 1.7
$ws.Range("T14").Value = 6
$ws.Range("U14").Value = 9.75
$ws.Range("V14").Value = 9.75
$ws.Range("W14").Value = 23
$ws.Range("X14").Value = 23
$ws.Range("Y14").Value = 40
$ws.Range("Z14").Value = 5.6
$ws.Range("AA14").Value = 5.9
$ws.Range("AB14").Value = 17.5
$ws.Range("AC14").Value = 110
$ws.Range("AE14").Value = 7.5
$ws.Range("AF14").Value = 15
$ws.Range("AG14").Value = 12
$ws.Range("AH14").Value = 45
$ws.Range("AI14").Value = 35
$ws.Range("AJ14").Value = 50
# Row 17
$ws.Range("G17").Value = 1.5
$ws.Range("H17").Value = 3.75
$ws.Range("I17").Value = 5.8
$ws.Range("N17").Value = 1.88
$ws.Range("O17").Value = 1.72
$ws.Range("T17").Value = 5.1
$ws.Range("U17").Value = 5.5
$ws.Range("V17").Value = 7
$ws.Range("X17").Value = 10.5
$ws.Range("Y17").Value = 24
$ws.Range("Z17").Value = 9.25
$ws.Range("AA17").Value = 6.5
$ws.Range("AB17").Value = 16
$ws.Range("AC17").Value = 75
$ws.Range("AE17").Value = 11.75
$ws.Range("AF17").Value = 28
$ws.Range("AG17").Value = 15.5
$ws.Range("AH17").Value = 90
$ws.Range("AI17").Value = 50
# Row 18
$ws.Range("G18").Value = 1.75
$ws.Range("H18").Value = 3.5
$ws.Range("I18").Value = 3.95
$ws.Range("N18").Value = 1.87
$ws.Range("O18").Value = 1.75
$ws.Range("T18").Value = 5.8
$ws.Range("U18").Value = 6.8
$ws.Range("V18").Value = 7.1
$ws.Range("X18").Value = 11.75
$ws.Range("Y18").Value = 22
$ws.Range("Z18").Value = 9.75
$ws.Range("AA18").Value = 6
$ws.Range("AB18").Value = 13.5
$ws.Range("AC18").Value = 60
$ws.Range("AD18").Value = 400
$ws.Range("AE18").Value = 9.25
$ws.Range("AF18").Value = 17
$ws.Range("AG18").Value = 11.25
$ws.Range("AH18").Value = 45
$ws.Range("AI18").Value = 30
$ws.Range("AJ18").Value = 35
# Row 25
$ws.Range("K25").Value = 9
# Row 28
$ws.Range("G28").Value = 1.44
$ws.Range("I28").Value = 6
$ws.Range("R28").Value = 1.83
$ws.Range("S28").Value = 1.83
$ws.Range("T28").Value = 7.5
$ws.Range("AF28").Value = 34
$ws.Range("AG28").Value = 19
# Row 30
$ws.Range("K30").Value = 13
# Row 31
$ws.Range("G31").Value = 1.8
$ws.Range("H31").Value = 3.65
$ws.Range("I31").Value = 3.9
$ws.Range("L31").Value = 1.24
$ws.Range("M31").Value = 3.3
$ws.Range("N31").Value = 1.72
$ws.Range("O31").Value = 1.88
$ws.Range("R31").Value = 1.7
$ws.Range("S31").Value = 1.93
$ws.Range("T31").Value = 7.8
$ws.Range("U31").Value = 9
$ws.Range("W31").Value = 15
$ws.Range("X31").Value = 13.5
$ws.Range("Y31").Value = 24
$ws.Range("Z31").Value = 11.5
$ws.Range("AA31").Value = 7.1
$ws.Range("AB31").Value = 14.5
$ws.Range("AC31").Value = 60
$ws.Range("AD31").Value = 450
$ws.Range("AE31").Value = 12
$ws.Range("AF31").Value = 22
$ws.Range("AJ31").Value = 40
# Row 32
$ws.Range("H32").Value = 4.05
$ws.Range("I32").Value = 4.35
$ws.Range("J32").Value = 1.04
$ws.Range("K32").Value = 8.75
$ws.Range("L32").Value = 1.2
$ws.Range("M32").Value = 4.05
$ws.Range("N32").Value = 1.62
$ws.Range("O32").Value = 2.18
$ws.Range("P32").Value = 1.31
$ws.Range("Q32").Value = 3.15
$ws.Range("R32").Value = 1.65
$ws.Range("S32").Value = 2.1
$ws.Range("T32").Value = 8.75
$ws.Range("U32").Value = 8.75
$ws.Range("X32").Value = 12
$ws.Range("Y32").Value = 21
$ws.Range("Z32").Value = 8.75
$ws.Range("AA32").Value = 8
$ws.Range("AB32").Value = 14.5
$ws.Range("AC32").Value = 55
$ws.Range("AD32").Value = 350
$ws.Range("AE32").Value = 15
$ws.Range("AF32").Value = 26
$ws.Range("AI32").Value = 37
$ws.Range("AJ32").Value = 37
# Row 33
$ws.Range("K33").Value = 13
